# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K") on Sheet1 held "Strike#" values; replace them with the
# correct strikeout (K) counts for each of the 37 data rows (rows 2-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 6
    3  = 2
    4  = 6
    5  = 5
    6  = 9
    7  = 9
    8  = 7
    9  = 3
    10 = 5
    11 = 12
    12 = 5
    13 = 5
    14 = 5
    15 = 5
    16 = 6
    17 = 7
    18 = 4
    19 = 5
    20 = 3
    21 = 2
    22 = 6
    23 = 4
    24 = 6
    25 = 9
    26 = 4
    27 = 10
    28 = 6
    29 = 6
    30 = 8
    31 = 7
    32 = 6
    33 = 7
    34 = 3
    35 = 5
    36 = 6
    37 = 2
    38 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
